$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are numeric-looking text that must be forced to
# remain text (otherwise Excel auto-converts them to numbers and drops
# significant trailing zeros, e.g. "69.60" -> 69.6).
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range('D2').Value = '51.605.24'
$ws.Range('E2').Value = '  +4.16%  '

$ws.Range('D3').Value = '2.759.88'
$ws.Range('E3').Value = '  +4.58%  '

$ws.Range('E4').Value = '  +0.09%  '

$ws.Range('D5').Value = '116.13'
$ws.Range('E5').Value = '  +3.21%  '

$ws.Range('D6').Value = '333.66'
$ws.Range('E6').Value = '  +2.99%  '

$ws.Range('D7').Value = '0.538'
$ws.Range('E7').Value = '  +2.07%  '

$ws.Range('E8').Value = '  +0.04%  '

$ws.Range('D9').Value = '0.575'
$ws.Range('E9').Value = '  +5.60%  '

$ws.Range('D10').Value = '41.82'
$ws.Range('E10').Value = '  +4.76%  '

$ws.Range('D11').Value = '0.0863'
$ws.Range('E11').Value = '  +6.10%  '

$ws.Range('D12').Value = '20.23'
$ws.Range('E12').Value = '  +1.98%  '

$ws.Range('E13').Value = '  +1.93%  '

$ws.Range('D14').Value = '7.65'
$ws.Range('E14').Value = '  +4.61%  '

$ws.Range('D15').Value = '3.194.69'
$ws.Range('E15').Value = '  +5.11%  '

$ws.Range('D16').Value = '2.765.88'
$ws.Range('E16').Value = '  +4.48%  '

$ws.Range('E17').Value = '  +3.32%  '

$ws.Range('D18').Value = '51.588.81'
$ws.Range('E18').Value = '  +4.44%  '

$ws.Range('D19').Value = '3.24'
$ws.Range('E19').Value = '  +7.44%  '

$ws.Range('D20').Value = '13.48'
$ws.Range('E20').Value = '  +4.34%  '

$ws.Range('E21').Value = '  +2.34%  '

$ws.Range('D22').Value = '0.0₃0973'
$ws.Range('E22').Value = '  +2.68%  '

$ws.Range('D23').Value = '278.51'
$ws.Range('E23').Value = '  +3.00%  '

$ws.Range('D24').Value = '69.60'
$ws.Range('E24').Value = '  +1.02%  '

$ws.Range('D25').Value = '2.67'
$ws.Range('E25').Value = '  +5.47%  '

$ws.Range('D26').Value = '26.85'
$ws.Range('E26').Value = '  +2.04%  '

$ws.Range('E27').Value = '  -0.01%  '

$ws.Range('D28').Value = '10.17'
$ws.Range('E28').Value = '  -1.40%  '

$ws.Range('E29').Value = '  +0.27%  '

$ws.Range('E30').Value = '  +1.96%  '

$ws.Range('D31').Value = '35.01'
$ws.Range('E31').Value = '  -0.40%  '

$ws.Range('D32').Value = '50.00'
$ws.Range('E32').Value = '  +0.83%  '

$ws.Range('E33').Value = '  +1.37%  '

$ws.Range('D34').Value = '0.0824'
$ws.Range('E34').Value = '  +1.37%  '

$ws.Range('E35').Value = '  +0.17%  '

$ws.Range('D36').Value = '18.97'
$ws.Range('E36').Value = '  +0.01%  '

$ws.Range('D37').Value = '4.99'
$ws.Range('E37').Value = '  +1.02%  '

$ws.Range('E38').Value = '  +1.57%  '

$ws.Range('E39').Value = '  +3.14%  '

$ws.Range('D40').Value = '0.0354'
$ws.Range('E40').Value = '  +9.12%  '

$ws.Range('D41').Value = '127.35'
$ws.Range('E41').Value = '  +0.69%  '

$ws.Range('D42').Value = '23.17'
$ws.Range('E42').Value = '  +3.57%  '

$ws.Range('E43').Value = '  +3.05%  '

$ws.Range('D44').Value = '2.30'
$ws.Range('E44').Value = '  +7.45%  '

$ws.Range('E45').Value = '  +14.67%  '

$ws.Range('D46').Value = '2.091.24'
$ws.Range('E46').Value = '  +1.47%  '

$ws.Range('E47').Value = '  +2.70%  '

$ws.Range('E48').Value = '  +4.58%  '

$ws.Range('E49').Value = '  +5.93%  '

$ws.Range('D50').Value = '8.99'
$ws.Range('E50').Value = '  +0.65%  '

$ws.Range('D51').Value = '59.90'
$ws.Range('E51').Value = '  +1.21%  '
